$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.07628696022732853
$ws.Range("H2").Value = 15.03234444479457
$ws.Range("I2").Value = -12.03444491373377

$ws.Range("G3").Value = 0.08504947791364285
$ws.Range("H3").Value = -28.08438665779851

$ws.Range("G4").Value = -0.3151499939134525
$ws.Range("H4").Value = -14.60193522882694

$ws.Range("G5").Value = -0.3282556134615772
$ws.Range("H5").Value = 17.73348734733268

$ws.Range("G6").Value = 0.1974506857778885
$ws.Range("H6").Value = 0.1529320065889339

$ws.Range("G7").Value = 0.3014597271804503
$ws.Range("H7").Value = 45.36426530507575

$ws.Range("G8").Value = 0.1234611316939594
$ws.Range("H8").Value = 21.16384416135609

$ws.Range("G9").Value = 0.154402335612279
$ws.Range("H9").Value = 22.07642382555814

$ws.Range("G10").Value = 0.03769039930365871
$ws.Range("H10").Value = -38.65275709310126

$ws.Range("G11").Value = 0.01916269804402507
$ws.Range("H11").Value = -61.62098927476305

$ws.Range("G12").Value = 0.070459616644268
$ws.Range("H12").Value = -23.88274888713201

$ws.Range("G13").Value = 0.1165999937804763
$ws.Range("H13").Value = 53.00286615009556

$ws.Range("G14").Value = 0.2437620798510717
$ws.Range("H14").Value = 7.877636443924155

$ws.Range("G15").Value = 0.2488479412876118
$ws.Range("H15").Value = 1.008828725353352

$ws.Range("G16").Value = 0.0971190274523487
$ws.Range("H16").Value = -14.61677548059599

$ws.Range("G17").Value = 0.1533093946331618
$ws.Range("H17").Value = 2.610971056441752

$ws.Range("G18").Value = -0.005837057972220477
$ws.Range("H18").Value = 34.79567558858483

$ws.Range("G19").Value = 0.02916049217734911
$ws.Range("H19").Value = 20.3895051595287

$ws.Range("G20").Value = 0.1301913552171849
$ws.Range("H20").Value = 53.0561728329204

$ws.Range("G21").Value = 0.0747842873191851
$ws.Range("H21").Value = 14.25478341812668

$ws.Range("G22").Value = 0.1868864693653663
$ws.Range("H22").Value = -2.440576529950394

$ws.Range("G23").Value = 0.1985637956798045
$ws.Range("H23").Value = -7.947587024195888

$ws.Range("G24").Value = -0.02238928487075818
$ws.Range("H24").Value = -488.3295514002349

$ws.Range("G25").Value = -0.01100491374560266
$ws.Range("H25").Value = 52.68193770877892

$ws.Range("G26").Value = 0.1892841667878095
$ws.Range("H26").Value = -7.606010818755013

$ws.Range("G27").Value = 0.2290079342876671
$ws.Range("H27").Value = 18.7279881292008

$ws.Range("G28").Value = 0.0652080022873248
$ws.Range("H28").Value = -2.549061711566954

$ws.Range("G29").Value = 0.06922626576626248
$ws.Range("H29").Value = -26.56173978128606
